# Generate Report for Handoff
# Updates the "Latest Handoff"/"Latest HO Xliff Generate Date" timestamps
# for the a9bf8ccf-... item (row 7) across the Overview, zh-cn and de-de
# sheets, reflecting a freshly generated handoff report.

$wb = $excel.ActiveWorkbook

# Overview sheet: column G = "Latest HO Xliff Generate Date" for row 7
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G7").Value = "2016-09-05 08:56:52"

# zh-cn sheet: column H = "Latest Handoff Datetime" for row 7
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H7").Value = "2016-09-05 08:56:47"

# de-de sheet: column H = "Latest Handoff Datetime" for row 7
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H7").Value = "2016-09-05 08:56:52"
